$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "metadata" sheet has a two-column form (A = field label, B = value) with
# guidance notes in column C. The contributor inserted a blank row above the
# "portrait" field (shifting the "portrait" / "description-bio" / "extra
# message" rows in columns A:B down by one row) while leaving the column C
# notes exactly where they were.

# Capture the values + bold formatting of the cells that are moving before
# overwriting anything.
$portraitLabel = $ws.Range("A8").Value2
$portraitBold = $ws.Range("A8").Font.Bold

$descLabel = $ws.Range("A9").Value2
$descBold = $ws.Range("A9").Font.Bold
$bioValue = $ws.Range("B9").Value2
$bioBold = $ws.Range("B9").Font.Bold

$extraLabel = $ws.Range("A10").Value2
$extraBold = $ws.Range("A10").Font.Bold

# Clear the old column A/B cells in rows 8-10 (their content is being moved,
# and the vacated cells should not leave behind any formatting residue).
$ws.Range("A8:B10").Clear()

# Write the moved content back one row lower, preserving bold formatting.
$ws.Range("A9").Value = $portraitLabel
$ws.Range("A9").Font.Bold = $portraitBold

$ws.Range("A10").Value = $descLabel
$ws.Range("A10").Font.Bold = $descBold
$ws.Range("B10").Value = $bioValue
$ws.Range("B10").Font.Bold = $bioBold

$ws.Range("A11").Value = $extraLabel
$ws.Range("A11").Font.Bold = $extraBold

# Update the active selection to match the new cursor position.
$null = $ws.Range("J7").Select()
